$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the last GPS log entry (Duren) with the new one (unna)
$ws.Range("A2").Value = "unna"
$ws.Range("B2").Value = "51.5333,7.6833"

# Shrink column B to fit the new, much shorter coordinate text
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Move the active selection to A2
$ws.Range("A2").Select() | Out-Null
